# Add the new game-results row (row 45) to the "data" worksheet and
# update the active selection to match (H46), per the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$row = 45

$ws.Cells.Item($row, 1).Value = 2
$ws.Cells.Item($row, 2).Value = "Alien Brood Encounters"
$ws.Cells.Item($row, 3).Value = "Mephisto"
$ws.Cells.Item($row, 4).Value = "Underworld|Aspects of the Void"
$ws.Cells.Item($row, 5).Value = "The Brood|Cape-Killers"
$ws.Cells.Item($row, 6).Value = "Dr. Octopus (V)|Captain America (B)|Namora (WW)|Caiera (WW)|Invisible Woman (FF)"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 11).Value = "Lots of wounds, but also lots of ways to KO them (or the twists before they spawn wounds). Did get to 5 escaped villains."
$ws.Cells.Item($row, 8).Value = "17|33"
$ws.Cells.Item($row, 9).Value = "not really"

$ws.Range("H46").Select()
